# review-service.xlsx seed-data update
#
# The review-service-data sheet stores one row per review, with column D
# ("ProductName") holding a shared product-name string that is repeated
# across every review row for that product. This change trims/fixes a
# batch of those product names (shortening overly-specific model numbers,
# fixing a typo, tweaking wording) and, in a few spots, introduces a
# slightly different variant of the name for a handful of rows within a
# product group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sony KLV-32R2027 -> Sony KLV   (rows 2-11)
for ($r = 2; $r -le 11; $r++) { $ws.Cells.Item($r, 4).Value = 'Sony KLV' }

# Samsung LED samrt TV -> Samsung LED TV   (rows 12-21, fixes "samrt" typo)
for ($r = 12; $r -le 21; $r++) { $ws.Cells.Item($r, 4).Value = 'Samsung LED TV' }

# LG LED smart TV -> LG LED  TV / LG LED TV   (rows 22-31)
# Row 22 keeps the (now double-spaced) original wording; rows 23-31 pick up
# the newly introduced "LG LED TV" variant.
$ws.Range("D22").Value = 'LG LED  TV'
for ($r = 23; $r -le 31; $r++) { $ws.Cells.Item($r, 4).Value = 'LG LED TV' }

# Canon EOS 1300D -> Canon EOS  / Canon EOS   (rows 32-36)
# Row 35 picks up the newly introduced "Canon EOS" (no trailing space) variant;
# the rest keep the trimmed "Canon EOS " (trailing space) wording.
for ($r = 32; $r -le 34; $r++) { $ws.Cells.Item($r, 4).Value = 'Canon EOS ' }
$ws.Range("D35").Value = 'Canon EOS'
$ws.Range("D36").Value = 'Canon EOS '

# Nikon D5600 Digital -> Nikon D5600    (rows 37-41)
for ($r = 37; $r -le 41; $r++) { $ws.Cells.Item($r, 4).Value = 'Nikon D5600 ' }

# Whirpool 0.8 ton 3 star Inverter split AC -> Whirpool split AC   (rows 42-46)
for ($r = 42; $r -le 46; $r++) { $ws.Cells.Item($r, 4).Value = 'Whirpool split AC' }

# Sanyo 1.5 3star Inverter split AC -> Sanyo split AC   (rows 47-51)
for ($r = 47; $r -le 51; $r++) { $ws.Cells.Item($r, 4).Value = 'Sanyo split AC' }

# Godrej 190 L3star Single door -> Godrej L3star   (rows 52-56)
for ($r = 52; $r -le 56; $r++) { $ws.Cells.Item($r, 4).Value = 'Godrej L3star' }

# LG 260 L 4Star Frost Free Double Door -> LG 4Star / LG 4Star    (rows 57-61)
# Rows 57-58 keep the trimmed "LG 4Star" wording; rows 59-61 pick up the
# newly introduced "LG 4Star " (trailing space) variant.
for ($r = 57; $r -le 58; $r++) { $ws.Cells.Item($r, 4).Value = 'LG 4Star' }
for ($r = 59; $r -le 61; $r++) { $ws.Cells.Item($r, 4).Value = 'LG 4Star ' }

# Restore the view: scrolled so column D / row 2 is the top-left visible
# cell, with D61 as the active selection.
try { $excel.ActiveWindow.ScrollColumn = 4 } catch {}
try { $excel.ActiveWindow.ScrollRow = 2 } catch {}
$ws.Range("D61").Select()
